# Commit: "Session 24 added - Spring Boot Security ACL"
#
# The "Session 21" block (row 105) is followed by a small bullet list in
# column B (AOP / Security - Permission Management on the WEB / Microservice /
# Cloud / Reactive / Kafka MQ / Spring Batch) that lives in rows 106-112,
# sharing row 106 with the "Security" / "Role Based Access Control" entry in
# columns C/D.
#
# This change inserts a brand new "Session 24" (Spring Security ACL) block,
# followed by an empty placeholder "Session 25" block, right after row 106 -
# pushing the old column-B bullet list further down the sheet (with a small
# gap of blank rows left above it, matching how Excel naturally leaves the
# previously-existing blank separator row in place as content shifts down).
# The "Microservice" bullet is also reworded to "Microservice, docker,
# kubernetes" while it is being moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The column-B bullet list currently anchored to row 106 (shared with the
# "Security" / "Role Based Access Control" row) needs to move out from under
# row 106 before we carve out space for the new session block.
$ws.Range("B106").ClearContents()

# Make room for the new Session 24 / Session 25 block: 12 blank rows inserted
# right after row 106 (6 rows for the new block + its trailing blank
# separator, 6 rows that the old bullet-list rows 107-112 shift into).
$ws.Range("A107:A118").EntireRow.Insert()

# --- New "Session 24" block -------------------------------------------------
$ws.Range("D108").Value = "Spring Security ACL"
$ws.Range("D110").Value = "https://github.com/spring-projects/spring-security/blob/main/acl/src/main/java/org/springframework/security/acls/domain/BasePermission.java"
$ws.Range("D109").Value = "https://www.baeldung.com/spring-security-acl"
$ws.Range("B120").Value = "Microservice, docker, kubernetes"
$ws.Range("D111").Value = "https://www.youtube.com/watch?v=GTln3jc5_eg"
$ws.Range("A107").Value = "Session 24"

# --- New empty "Session 25" placeholder block -------------------------------
$ws.Range("A112").Value = "Session 25"

# Give the two new section-header rows the same look (bold/filled, merged
# A:E) as every other session header row on the sheet.
$ws.Range("A105:E105").Copy()
$ws.Range("A107:E107").PasteSpecial(-4122)
$ws.Range("A112:E112").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A107:E107").Merge()
$ws.Range("A112:E112").Merge()

# The old column-B bullet list (now starting a few rows further down) keeps
# its first item, "AOP", which was cleared off of row 106 above.
$ws.Range("B118").Value = "AOP"

$ws.Range("D113").Select()
